# Update elapsed-duration ("Elapsed Duration(Hrs)") readings and fill in the
# newly-reported outage row on sheet R1, per the Active_Outages.xlsx refresh.

$wb = $excel.ActiveWorkbook

# --- R1 ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3918:36:01"
$ws1.Range("G3").Value = "58:08:39"

# Newly populated outage record in row 5
$ws1.Range("B5").Value = "R4"
$ws1.Range("D5").Value = "JED0123"
$ws1.Range("I5").Value = "SCECO"
$ws1.Range("J5").Value = "In progress"
$ws1.Range("L5").Value = "Latis"

# --- R2 ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12099:59:42"
$ws2.Range("G3").Value = "3229:43:11"
$ws2.Range("G4").Value = "467:54:45"

# --- R4 ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2945:49:31"
$ws4.Range("G3").Value = "173:01:46"

# --- R5 ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "419:48:30"

# --- R6 ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "60:20:48"
